$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report period / volume number) ---
$ws.Range("A8").Value = "Volume 32   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/27/2025  Through  2/2/2025"

# --- Cells changing type (numeric <-> text "N/A"/"***.*") ---
# Use Copy() from a same-style/type neighbor cell so the destination
# picks up the correct style + shared-string text type, then overwrite
# the value for cells that become numeric.
$ws.Range("C14").Copy($ws.Range("D14"))     # D14: 2 -> "0" (N/A)
$ws.Range("M14").Copy($ws.Range("E14"))     # E14: -100 -> "***.*"
$ws.Range("D15").Copy($ws.Range("C15"))     # C15: 1 -> "0" (N/A)
$ws.Range("I22").Copy($ws.Range("C18"))     # C18: "0" (N/A) -> 1
$ws.Range("C18").Value = 1
$ws.Range("D27").Copy($ws.Range("C27"))     # C27: 1 -> "0" (N/A)
$ws.Range("C29").Copy($ws.Range("D29"))     # D29: 2 -> "0" (N/A)
$ws.Range("E15").Copy($ws.Range("E29"))     # E29: -100 -> "***.*"
$ws.Range("C30").Copy($ws.Range("D30"))     # D30: 2 -> "0" (N/A)
$ws.Range("H15").Copy($ws.Range("E30"))     # E30: -100 -> "***.*"

# --- Remaining numeric value updates ---
    $ws.Range("M15").Value = -50
    $ws.Range("N15").Value = -40
    $ws.Range("C16").Value = 7
    $ws.Range("D16").Value = 9
    $ws.Range("E16").Value = -22.222222222222
    $ws.Range("F16").Value = 20
    $ws.Range("G16").Value = 26
    $ws.Range("H16").Value = -23.076923076923
    $ws.Range("I16").Value = 21
    $ws.Range("J16").Value = 29
    $ws.Range("K16").Value = -27.586206896551
    $ws.Range("L16").Value = -8.695652173913
    $ws.Range("M16").Value = -34.375
    $ws.Range("N16").Value = -89.447236180904
    $ws.Range("C17").Value = 7
    $ws.Range("D17").Value = 8
    $ws.Range("E17").Value = -12.5
    $ws.Range("F17").Value = 39
    $ws.Range("G17").Value = 46
    $ws.Range("H17").Value = -15.217391304347
    $ws.Range("I17").Value = 49
    $ws.Range("J17").Value = 56
    $ws.Range("L17").Value = -26.865671641791
    $ws.Range("N17").Value = -49.484536082474
    $ws.Range("D18").Value = 2
    $ws.Range("E18").Value = -50
    $ws.Range("F18").Value = 4
    $ws.Range("G18").Value = 7
    $ws.Range("H18").Value = -42.857142857142
    $ws.Range("I18").Value = 6
    $ws.Range("J18").Value = 9
    $ws.Range("K18").Value = -33.333333333333
    $ws.Range("L18").Value = -70
    $ws.Range("M18").Value = -76.923076923076
    $ws.Range("C19").Value = 1
    $ws.Range("D19").Value = 12
    $ws.Range("E19").Value = -91.666666666666
    $ws.Range("F19").Value = 13
    $ws.Range("G19").Value = 32
    $ws.Range("H19").Value = -59.375
    $ws.Range("I19").Value = 17
    $ws.Range("J19").Value = 39
    $ws.Range("K19").Value = -56.410256410256
    $ws.Range("L19").Value = -54.054054054054
    $ws.Range("M19").Value = -52.777777777777
    $ws.Range("N19").Value = -70.689655172413
    $ws.Range("C20").Value = 3
    $ws.Range("D20").Value = 2
    $ws.Range("E20").Value = 50
    $ws.Range("F20").Value = 8
    $ws.Range("G20").Value = 13
    $ws.Range("H20").Value = -38.461538461538
    $ws.Range("I20").Value = 10
    $ws.Range("J20").Value = 17
    $ws.Range("K20").Value = -41.176470588235
    $ws.Range("L20").Value = -16.666666666666
    $ws.Range("M20").Value = 0
    $ws.Range("N20").Value = -83.333333333333
    $ws.Range("C21").Value = 19
    $ws.Range("D21").Value = 33
    $ws.Range("E21").Value = -42.424242424242
    $ws.Range("F21").Value = 87
    $ws.Range("G21").Value = 127
    $ws.Range("H21").Value = -31.496062992126
    $ws.Range("I21").Value = 106
    $ws.Range("J21").Value = 153
    $ws.Range("K21").Value = -30.718954248366
    $ws.Range("L21").Value = -35.365853658536
    $ws.Range("M21").Value = -30.263157894736
    $ws.Range("N21").Value = -78.189300411522
    $ws.Range("D22").Value = 2
    $ws.Range("F22").Value = 2
    $ws.Range("G22").Value = 11
    $ws.Range("H22").Value = -81.818181818181
    $ws.Range("J22").Value = 12
    $ws.Range("K22").Value = -75
    $ws.Range("L22").Value = -50
    $ws.Range("M22").Value = 50
    $ws.Range("C23").Value = 4
    $ws.Range("D23").Value = 3
    $ws.Range("E23").Value = 33.333333333333
    $ws.Range("F23").Value = 17
    $ws.Range("G23").Value = 23
    $ws.Range("H23").Value = -26.086956521739
    $ws.Range("I23").Value = 25
    $ws.Range("J23").Value = 28
    $ws.Range("K23").Value = -10.714285714285
    $ws.Range("L23").Value = -28.571428571428
    $ws.Range("M23").Value = 13.636363636363
    $ws.Range("C24").Value = 14
    $ws.Range("D24").Value = 27
    $ws.Range("E24").Value = -48.148148148148
    $ws.Range("F24").Value = 70
    $ws.Range("G24").Value = 92
    $ws.Range("H24").Value = -23.913043478260
    $ws.Range("I24").Value = 79
    $ws.Range("J24").Value = 104
    $ws.Range("K24").Value = -24.038461538461
    $ws.Range("L24").Value = -36.8
    $ws.Range("M24").Value = -16.842105263157
    $ws.Range("C25").Value = 1
    $ws.Range("D25").Value = 4
    $ws.Range("E25").Value = -75
    $ws.Range("F25").Value = 12
    $ws.Range("G25").Value = 27
    $ws.Range("H25").Value = -55.555555555555
    $ws.Range("I25").Value = 15
    $ws.Range("J25").Value = 30
    $ws.Range("K25").Value = -50
    $ws.Range("L25").Value = -59.459459459459
    $ws.Range("C26").Value = 16
    $ws.Range("D26").Value = 9
    $ws.Range("E26").Value = 77.777777777777
    $ws.Range("F26").Value = 63
    $ws.Range("G26").Value = 52
    $ws.Range("H26").Value = 21.153846153846
    $ws.Range("I26").Value = 76
    $ws.Range("J26").Value = 65
    $ws.Range("K26").Value = 16.923076923076
    $ws.Range("L26").Value = -16.483516483516
    $ws.Range("M26").Value = -28.301886792452
    $ws.Range("C28").Value = 2
    $ws.Range("D28").Value = 2
    $ws.Range("F28").Value = 7
    $ws.Range("G28").Value = 4
    $ws.Range("H28").Value = 75
    $ws.Range("I28").Value = 10
    $ws.Range("J28").Value = 7
    $ws.Range("K28").Value = 42.857142857142
    $ws.Range("L28").Value = -9.090909090909
    $ws.Range("G29").Value = 7
    $ws.Range("H29").Value = -85.714285714285
    $ws.Range("L29").Value = -83.333333333333
    $ws.Range("N29").Value = -95.833333333333
    $ws.Range("G30").Value = 6
    $ws.Range("H30").Value = -83.333333333333
    $ws.Range("L30").Value = -83.333333333333
    $ws.Range("N30").Value = -95.833333333333
